$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 5337
$ws.Range("E3").Value = 10711
$ws.Range("E4").Value = 14139
$ws.Range("E5").Value = 9677
$ws.Range("E6").Value = 4914
$ws.Range("E7").Value = 9715
$ws.Range("E8").Value = 9914
$ws.Range("E9").Value = 10246
$ws.Range("E10").Value = 13883
$ws.Range("E11").Value = 1920
$ws.Range("E12").Value = 13323
$ws.Range("E13").Value = 3046
